# Applies the cryptos-list price/volume refresh described by the commit
# "Updated cryptos list on Sat Apr 27 03:57:22 UTC 2024 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '63.128.60'
$ws.Range("E2").Value2 = '  -2.03%  '
$ws.Range("D3").Value2 = '3.129.59'
$ws.Range("E3").Value2 = '  -0.49%  '
$ws.Range("E4").Value2 = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '593.79'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '136.70'
$ws.Range("E6").Value2 = '  -4.80%  '
$ws.Range("D8").Value2 = '3.120.66'
$ws.Range("E8").Value2 = '  -0.74%  '
$ws.Range("E10").Value2 = '  -2.68%  '
$ws.Range("E11").Value2 = '  -0.53%  '
$ws.Range("E12").Value2 = '  -2.80%  '
$ws.Range("E13").Value2 = '  -2.58%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '34.22'
$ws.Range("E14").Value2 = '  -3.76%  '
$ws.Range("D15").Value2 = '3.641.05'
$ws.Range("E15").Value2 = '  -0.60%  '
$ws.Range("E16").Value2 = '  +1.19%  '
$ws.Range("D17").Value2 = '63.118.38'
$ws.Range("E17").Value2 = '  -1.96%  '
$ws.Range("D18").Value2 = '3.111.30'
$ws.Range("E18").Value2 = '  -1.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '6.73'
$ws.Range("E19").Value2 = '  -1.72%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '477.79'
$ws.Range("E20").Value2 = '  +0.13%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '14.21'
$ws.Range("E21").Value2 = '  -3.02%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '0.700'
$ws.Range("E22").Value2 = '  -3.29%  '
$ws.Range("E23").Value2 = '  -1.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '87.84'
$ws.Range("E24").Value2 = '  +3.43%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '13.13'
$ws.Range("E25").Value2 = '  -3.59%  '
$ws.Range("E26").Value2 = '  +0.18%  '
$ws.Range("E27").Value2 = '  -2.93%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '7.24'
$ws.Range("E28").Value2 = '  -2.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '8.01'
$ws.Range("E29").Value2 = '  -6.61%  '
$ws.Range("E30").Value2 = '  -0.19%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '27.23'
$ws.Range("E31").Value2 = '  +2.05%  '
$ws.Range("E32").Value2 = '  -0.01%  '
$ws.Range("E33").Value2 = '  -7.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '2.56'
$ws.Range("E34").Value2 = '  -3.26%  '
$ws.Range("E35").Value2 = '  -3.10%  '
$ws.Range("E36").Value2 = '  -1.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '52.07'
$ws.Range("E37").Value2 = '  -0.94%  '
$ws.Range("D38").Value2 = '0.0₃0716'
$ws.Range("E38").Value2 = '  -3.41%  '
$ws.Range("E39").Value2 = '  -1.20%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '421.53'
$ws.Range("E40").Value2 = '  -6.82%  '
$ws.Range("E41").Value2 = '  -1.37%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '8.27'
$ws.Range("E42").Value2 = '  -0.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '2.70'
$ws.Range("E43").Value2 = '  -10.37%  '
$ws.Range("D44").Value2 = '2.885.18'
$ws.Range("E44").Value2 = '  +1.29%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '0.261'
$ws.Range("E45").Value2 = '  -2.40%  '
$ws.Range("E46").Value2 = '  -0.11%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '2.13'
$ws.Range("E47").Value2 = '  -5.40%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '25.60'
$ws.Range("E48").Value2 = '  -2.87%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '0.114'
$ws.Range("E49").Value2 = '  -0.22%  '
$ws.Range("E50").Value2 = '  -5.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '119.21'
$ws.Range("E51").Value2 = '  -1.26%  '
